$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'290.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-3.70%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'30.77"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-6.27%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.952"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.13%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07207"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-7.43%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.787"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-9.18%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.673"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-2.29%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.755"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.15%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.8968"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-3.21%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1648"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-6.54%"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'-2.69%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08018"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-6.91%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03038"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-3.52%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.1003"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.02%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001496"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.19%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005713"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.91%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.460"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.05%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'-3.30%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3314"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.59%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1315"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.25%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.046"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-6.54%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2185"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'9.72%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'-1.11%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-0.86%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004013"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-9.62%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-0.04%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.01605"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-6.44%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04395"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-6.96%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007380"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-6.01%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1308"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-3.43%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007661"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.002007"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-14.30%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009212"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-12.17%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00005946"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-4.88%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.03%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.246"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'173.66%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003002"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-3.26%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.03%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'-0.03%"
$ws.Range("E51").Style = "Normal"
